# Update the cryptocurrency price / volume table with the latest scrape.
# Values in column D that look numeric are written with a leading
# apostrophe so Excel stores them verbatim as text (preserving exact
# formatting such as trailing zeros, thousands separators written with
# dots, and leading zeros) instead of re-interpreting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.975.57"
$ws.Range("D3").Value = "'2.643.91"
$ws.Range("E3").Value = "  +2.21%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'538.39"
$ws.Range("E5").Value = "  +2.12%  "
$ws.Range("D6").Value = "'144.81"
$ws.Range("E6").Value = "  +4.15%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("E9").Value = "  +2.55%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").Value = "'0.338"
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "'3.110.08"
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").Value = "'59.867.46"
$ws.Range("E14").Value = "  +1.57%  "
$ws.Range("D15").Value = "'20.99"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000135"
$ws.Range("E16").Value = "  +1.49%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.596.72"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'343.42"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("D19").Value = "'4.41"
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("D20").Value = "'10.24"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").Value = "'0.414"
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("E25").Value = "  -0.74%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "'7.29"
$ws.Range("E27").Value = "  +3.16%  "
$ws.Range("D28").Value = "'" + "0.0" + [char]8323 + "0754"
$ws.Range("E28").Value = "  +4.62%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +3.44%  "
$ws.Range("D31").Value = "'5.86"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("D32").Value = "'18.95"
$ws.Range("E32").Value = "  +1.22%  "
$ws.Range("D33").Value = "'150.79"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").Value = "'0.843"
$ws.Range("E37").Value = "  +1.87%  "
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").Value = "'291.61"
$ws.Range("E39").Value = "  +8.28%  "
$ws.Range("D40").Value = "'3.59"
$ws.Range("E40").Value = "  +1.76%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.606"
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("D43").Value = "'10.74"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +3.56%  "
$ws.Range("D46").Value = "'1.977.33"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0225"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'18.66"
$ws.Range("E48").Value = "  +2.43%  "
$ws.Range("D49").Value = "'4.56"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("D50").Value = "'111.90"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("D51").Value = "'4.73"
$ws.Range("E51").Value = "  +0.16%  "
